$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.067668914794922
$ws.Range("C2").Value = 5.896552085876465
$ws.Range("D2").Value = 14.725564002990723
$ws.Range("E2").Value = 57.85714340209961
